$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0

$ws.Range("B3").Value = 0
$ws.Range("C3").Value = -0.6989857865567108
$ws.Range("D3").Value = 0

$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 0.6484505050432346
$ws.Range("D4").Value = 0

$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0

$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0

$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 0.666829450509164
$ws.Range("D7").Value = 0

$ws.Range("B8").Value = 0
$ws.Range("C8").Value = -0.6960145273620598
$ws.Range("D8").Value = 0

$ws.Range("B9").Value = 0
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 0
